# Weekly "Updated cryptos list" refresh: new prices / 1h-volume deltas for each
# coin row, plus two rank-tie swaps (HuobiToken<->TrustWalletToken at rows 36-37,
# WEMIXToken<->mCoin at rows 43-44) where the name/link/price/volume moved together.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price strings (col D) need a text quote-prefix so Excel keeps them
# as literal strings ("211.10") instead of coercing to the number 211.1.
$quote = [string][char]39

$ws.Cells.Item(2, 4).Value = '27.450.73'
$ws.Cells.Item(2, 5).Value = '  -1.25%  '
$ws.Cells.Item(3, 4).Value = '1.615.38'
$ws.Cells.Item(3, 5).Value = '  -2.15%  '
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).Value = $quote + '211.10'
$ws.Cells.Item(5, 5).Value = '  -1.18%  '
$ws.Cells.Item(6, 5).Value = '  -1.44%  '
$ws.Cells.Item(7, 5).Value = '  +0.13%  '
$ws.Cells.Item(8, 4).Value = $quote + '22.77'
$ws.Cells.Item(8, 5).Value = '  -1.88%  '
$ws.Cells.Item(10, 5).Value = '  -0.74%  '
$ws.Cells.Item(11, 4).Value = $quote + '0.0886'
$ws.Cells.Item(11, 5).Value = '  -0.61%  '
$ws.Cells.Item(12, 4).Value = '1.845.48'
$ws.Cells.Item(12, 5).Value = '  -2.03%  '
$ws.Cells.Item(13, 4).Value = '1.614.12'
$ws.Cells.Item(13, 5).Value = '  -2.25%  '
$ws.Cells.Item(15, 5).Value = '  -3.10%  '
$ws.Cells.Item(16, 4).Value = $quote + '65.09'
$ws.Cells.Item(16, 5).Value = '  +0.81%  '
$ws.Cells.Item(17, 4).Value = '27.425.59'
$ws.Cells.Item(17, 5).Value = '  -1.23%  '
$ws.Cells.Item(18, 4).Value = $quote + '233.11'
$ws.Cells.Item(18, 5).Value = '  -0.60%  '
$ws.Cells.Item(19, 5).Value = '  -1.32%  '
$ws.Cells.Item(20, 4).Value = $quote + '7.51'
$ws.Cells.Item(20, 5).Value = '  -2.75%  '
$ws.Cells.Item(21, 5).Value = '  +0.12%  '
$ws.Cells.Item(23, 4).Value = $quote + '10.21'
$ws.Cells.Item(23, 5).Value = '  +0.47%  '
$ws.Cells.Item(24, 5).Value = '  +5.47%  '
$ws.Cells.Item(25, 4).Value = $quote + '150.22'
$ws.Cells.Item(25, 5).Value = '  -0.31%  '
$ws.Cells.Item(26, 5).Value = '  -1.93%  '
$ws.Cells.Item(27, 5).Value = '  -1.58%  '
$ws.Cells.Item(28, 5).Value = '  +0.10%  '
$ws.Cells.Item(29, 5).Value = '  -1.20%  '
$ws.Cells.Item(30, 4).Value = $quote + '1.17'
$ws.Cells.Item(30, 5).Value = '  -1.52%  '
$ws.Cells.Item(31, 5).Value = '  -1.30%  '
$ws.Cells.Item(32, 5).Value = '  -1.57%  '
$ws.Cells.Item(33, 4).Value = '1.467.91'
$ws.Cells.Item(33, 5).Value = '  +1.67%  '
$ws.Cells.Item(34, 5).Value = '  -3.46%  '
$ws.Cells.Item(35, 5).Value = '  -3.69%  '
$ws.Cells.Item(36, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(36, 4).Value = $quote + '0.967'
$ws.Cells.Item(36, 5).Value = '  +10.00%  '
$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(37, 4).Value = $quote + '2.33'
$ws.Cells.Item(37, 5).Value = '  -0.50%  '
$ws.Cells.Item(38, 5).Value = '  -0.74%  '
$ws.Cells.Item(39, 5).Value = '  -2.73%  '
$ws.Cells.Item(40, 5).Value = '  -3.13%  '
$ws.Cells.Item(41, 5).Value = '  +0.10%  '
$ws.Cells.Item(42, 4).Value = $quote + '66.79'
$ws.Cells.Item(42, 5).Value = '  +0.09%  '
$ws.Cells.Item(43, 2).Value = 'mCoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Cells.Item(43, 4).Value = $quote + '2.46'
$ws.Cells.Item(43, 5).Value = '  -0.28%  '
$ws.Cells.Item(44, 2).Value = 'WEMIXToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(44, 4).Value = $quote + '0.984'
$ws.Cells.Item(44, 5).Value = '  -4.90%  '
$ws.Cells.Item(45, 5).Value = '  -2.82%  '
$ws.Cells.Item(46, 4).Value = $quote + '5.24'
$ws.Cells.Item(46, 5).Value = '  -6.79%  '
$ws.Cells.Item(47, 4).Value = '1.755.54'
$ws.Cells.Item(47, 5).Value = '  -2.09%  '
$ws.Cells.Item(48, 5).Value = '  -0.85%  '
$ws.Cells.Item(49, 4).Value = $quote + '86.77'
$ws.Cells.Item(49, 5).Value = '  +0.31%  '
$ws.Cells.Item(50, 4).Value = '0.0₆0105'
$ws.Cells.Item(51, 5).Value = '  +0.69%  '
